# Edit: Sun, May 10, 2020 8:19:00 AM
#
# 1) Slide 6's table switches to a different (built-in) table style.
# 2) The deck's applied theme ("Integral") swaps its colour scheme with the
#    previously-unused default "Office Theme" colour scheme, so the design
#    that actually renders the slides now uses the stock Office palette
#    instead of the green/gold Integral palette.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------
$s = $p.Slides.Item(6)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{C7B9EF85-47F4-4742-A725-BB1B81E18FB5}")

# --- 2) Swap the theme colour scheme (Integral -> Office Theme) ----------
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

# dk1 / lt1 / dk2 / lt2
$cs.Item(1).RGB  = 0          # dk1  000000
$cs.Item(2).RGB  = 16777215   # lt1  FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2  44546A
$cs.Item(4).RGB  = 15132391   # lt2  E7E6E6

# accent1-6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47

# hyperlink / followed hyperlink
$cs.Item(11).RGB = 12673797   # hlink     0563C1
$cs.Item(12).RGB = 7491477    # folHlink  954F72
